# Auto-generated: refresh market-price-derived columns (H..N) across all class sheets
# per scheduled-runner data update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1036.9642
$ws.Range("I17").Value = 385.7143
$ws.Range("J17").Value = 1130
$ws.Range("K17").Value = 1157.1429
$ws.Range("L17").Value = 3390
$ws.Range("M17").Value = -989.1428999999998
$ws.Range("N17").Value = -3726
$ws.Range("H69").Value = 9498
$ws.Range("I69").Value = 9498
$ws.Range("K69").Value = 28494
$ws.Range("M69").Value = -27620
$ws.Range("H72").Value = 9498
$ws.Range("I72").Value = 9498
$ws.Range("K72").Value = 85482
$ws.Range("M72").Value = -81114
$ws.Range("H96").Value = 998.5
$ws.Range("J96").Value = 998.5
$ws.Range("L96").Value = 2995.5
$ws.Range("N96").Value = -5741.5
$ws.Range("H98").Value = 7731.1113
$ws.Range("I98").Value = 7731.1113
$ws.Range("K98").Value = 7731.1113
$ws.Range("M98").Value = -6233.1113
$ws.Range("H116").Value = 27786056
$ws.Range("I116").Value = 50005200
$ws.Range("J116").Value = 12124.5
$ws.Range("K116").Value = 50005200
$ws.Range("L116").Value = 12124.5
$ws.Range("M116").Value = -50001758
$ws.Range("N116").Value = -19008.5
$ws.Range("H122").Value = 7731.1113
$ws.Range("I122").Value = 7731.1113
$ws.Range("K122").Value = 23193.3339
$ws.Range("M122").Value = -20743.3339
$ws.Range("H131").Value = 2491.353
$ws.Range("I131").Value = 1873.4615
$ws.Range("K131").Value = 5620.3845
$ws.Range("M131").Value = -580.3845000000001
$ws.Range("H135").Value = 625830.1
$ws.Range("I135").Value = 667472.1
$ws.Range("K135").Value = 6007248.899999999
$ws.Range("M135").Value = -6004713.899999999
$ws.Range("H137").Value = 2335.6428
$ws.Range("I137").Value = 2419
$ws.Range("K137").Value = 7257
$ws.Range("M137").Value = -4707
$ws.Range("H138").Value = 5822
$ws.Range("J138").Value = 9999.75
$ws.Range("L138").Value = 29999.25
$ws.Range("N138").Value = -40279.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3072.7
$ws.Range("I2").Value = 1371.6
$ws.Range("J2").Value = 4773.8
$ws.Range("K2").Value = 1371.6
$ws.Range("L2").Value = 4773.8
$ws.Range("M2").Value = -1258.6
$ws.Range("N2").Value = -4999.8
$ws.Range("H116").Value = 3072.7
$ws.Range("I116").Value = 1371.6
$ws.Range("J116").Value = 4773.8
$ws.Range("K116").Value = 1371.6
$ws.Range("L116").Value = 4773.8
$ws.Range("M116").Value = 922.4000000000001
$ws.Range("N116").Value = -9361.799999999999
$ws.Range("H132").Value = 4833.0493
$ws.Range("I132").Value = 3700.4888
$ws.Range("J132").Value = 8018.375
$ws.Range("K132").Value = 11101.4664
$ws.Range("L132").Value = 24055.125
$ws.Range("M132").Value = -8571.466400000001
$ws.Range("N132").Value = -29115.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3072.7
$ws.Range("I3").Value = 1371.6
$ws.Range("J3").Value = 4773.8
$ws.Range("K3").Value = 1371.6
$ws.Range("L3").Value = 4773.8
$ws.Range("M3").Value = -1257.6
$ws.Range("N3").Value = -5001.8
$ws.Range("H94").Value = 1190.25
$ws.Range("I94").Value = 636.1667
$ws.Range("K94").Value = 636.1667
$ws.Range("M94").Value = -185.1667
$ws.Range("H105").Value = 52057.516
$ws.Range("I105").Value = 66093.375
$ws.Range("K105").Value = 66093.375
$ws.Range("M105").Value = -64346.375
$ws.Range("H134").Value = 4372.339
$ws.Range("I134").Value = 1619.8462
$ws.Range("K134").Value = 4859.5386
$ws.Range("M134").Value = -2324.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7841.5684
$ws.Range("I31").Value = 2671.4546
$ws.Range("K31").Value = 2671.4546
$ws.Range("M31").Value = -2376.4546
$ws.Range("H34").Value = 7841.5684
$ws.Range("I34").Value = 2671.4546
$ws.Range("K34").Value = 2671.4546
$ws.Range("M34").Value = -2469.4546
$ws.Range("H68").Value = 84999.25
$ws.Range("J68").Value = 84999.25
$ws.Range("L68").Value = 84999.25
$ws.Range("N68").Value = -86497.25
$ws.Range("H71").Value = 84999.25
$ws.Range("J71").Value = 84999.25
$ws.Range("L71").Value = 254997.75
$ws.Range("N71").Value = -262485.75
$ws.Range("H86").Value = 10547000
$ws.Range("I86").Value = 20837834
$ws.Range("J86").Value = 256166.67
$ws.Range("K86").Value = 20837834
$ws.Range("L86").Value = 256166.67
$ws.Range("M86").Value = -20836711
$ws.Range("N86").Value = -258412.67
$ws.Range("H89").Value = 10547000
$ws.Range("I89").Value = 20837834
$ws.Range("J89").Value = 256166.67
$ws.Range("K89").Value = 104189170
$ws.Range("L89").Value = 1280833.35
$ws.Range("M89").Value = -104183554
$ws.Range("N89").Value = -1292065.35
$ws.Range("H132").Value = 4621.864
$ws.Range("I132").Value = 1512.1154
$ws.Range("J132").Value = 9113.723
$ws.Range("K132").Value = 4536.3462
$ws.Range("L132").Value = 27341.169
$ws.Range("M132").Value = -2006.3462
$ws.Range("N132").Value = -32401.169

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 78407016
$ws.Range("J4").Value = 800665.4399999999
$ws.Range("L4").Value = 2401996.32
$ws.Range("N4").Value = -2402220.32
$ws.Range("H117").Value = 528.1429000000001
$ws.Range("J117").Value = 599.6667
$ws.Range("L117").Value = 1799.0001
$ws.Range("N117").Value = -8683.000099999999
$ws.Range("H122").Value = 1770454.8
$ws.Range("I122").Value = 4042043.8
$ws.Range("J122").Value = 3663.2222
$ws.Range("K122").Value = 36378394.2
$ws.Range("L122").Value = 32968.99980000001
$ws.Range("M122").Value = -36375944.2
$ws.Range("N122").Value = -37868.99980000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9872.823
$ws.Range("I70").Value = 9134.5
$ws.Range("J70").Value = 10275.546
$ws.Range("K70").Value = 9134.5
$ws.Range("L70").Value = 10275.546
$ws.Range("M70").Value = -8864.5
$ws.Range("N70").Value = -10815.546
$ws.Range("H73").Value = 9872.823
$ws.Range("I73").Value = 9134.5
$ws.Range("J73").Value = 10275.546
$ws.Range("K73").Value = 9134.5
$ws.Range("L73").Value = 10275.546
$ws.Range("M73").Value = -8198.5
$ws.Range("N73").Value = -12147.546
$ws.Range("H102").Value = 3579.0667
$ws.Range("I102").Value = 3481.2144
$ws.Range("J102").Value = 4949
$ws.Range("K102").Value = 3481.2144
$ws.Range("L102").Value = 4949
$ws.Range("M102").Value = -1859.2144
$ws.Range("N102").Value = -8193
$ws.Range("H113").Value = 6088.2896
$ws.Range("J113").Value = 8333.261
$ws.Range("L113").Value = 8333.261
$ws.Range("N113").Value = -12673.261
$ws.Range("H132").Value = 6611.5557
$ws.Range("J132").Value = 13412.571
$ws.Range("L132").Value = 40237.713
$ws.Range("N132").Value = -45297.713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 37037416
$ws.Range("I55").Value = 90909160
$ws.Range("J55").Value = 594.8125
$ws.Range("K55").Value = 90909160
$ws.Range("L55").Value = 594.8125
$ws.Range("M55").Value = -90908987
$ws.Range("N55").Value = -940.8125
$ws.Range("H68").Value = 5499.75
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251
$ws.Range("H71").Value = 5499.75
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256
$ws.Range("H82").Value = 742641.5
$ws.Range("I82").Value = 1006928.94
$ws.Range("K82").Value = 1006928.94
$ws.Range("M82").Value = -1006567.94
$ws.Range("H85").Value = 742641.5
$ws.Range("I85").Value = 1006928.94
$ws.Range("K85").Value = 1006928.94
$ws.Range("M85").Value = -1005680.94
$ws.Range("H100").Value = 4588
$ws.Range("J100").Value = 6200.2
$ws.Range("L100").Value = 6200.2
$ws.Range("N100").Value = -7282.2
$ws.Range("H132").Value = 10422844
$ws.Range("I132").Value = 22729432
$ws.Range("J132").Value = 9576.654
$ws.Range("K132").Value = 68188296
$ws.Range("L132").Value = 28729.962
$ws.Range("M132").Value = -68185766
$ws.Range("N132").Value = -33789.962

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4656.8335
$ws.Range("I96").Value = 4656.8335
$ws.Range("K96").Value = 4656.8335
$ws.Range("M96").Value = -3283.8335
$ws.Range("H100").Value = 919.8823
$ws.Range("I100").Value = 693.8889
$ws.Range("J100").Value = 1174.125
$ws.Range("K100").Value = 1387.7778
$ws.Range("L100").Value = 2348.25
$ws.Range("M100").Value = -846.7778000000001
$ws.Range("N100").Value = -3430.25
$ws.Range("H136").Value = 37043630
$ws.Range("J136").Value = 8786
$ws.Range("L136").Value = 26358
$ws.Range("N136").Value = -31458

Write-Output "Updated 220 cells across 8 sheets"
